$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retardos3")

# Update R8:S17 values
$ws.Range("R8").Value = 5941
$ws.Range("S8").Value = 1694510

$ws.Range("R9").Value = 5611
$ws.Range("S9").Value = 1685108

$ws.Range("R10").Value = 6891
$ws.Range("S10").Value = 1696078

$ws.Range("R11").Value = 7662
$ws.Range("S11").Value = 1693462

$ws.Range("R12").Value = 6258
$ws.Range("S12").Value = 1697605

$ws.Range("R13").Value = 6897
$ws.Range("S13").Value = 1692796

$ws.Range("R14").Value = 5295
$ws.Range("S14").Value = 1693505

$ws.Range("R15").Value = 6584
$ws.Range("S15").Value = 1695937

$ws.Range("R16").Value = 6579
$ws.Range("S16").Value = 1692608

$ws.Range("R17").Value = 5622
$ws.Range("S17").Value = 1696778

# Add textbox to sheet3 (mirrors the "Codigo 1 sin CK" box already present on Retardos1/Retardos2)
$shp = $ws.Shapes.AddTextbox(1, 258.6, 2.4, 579.1207874, 52.8)
$shp.Name = "CuadroTexto 1"
$shp.Fill.ForeColor.RGB = 16777215
$line1 = "Código 1 sin CK"
$line2 = "Se envía una trama desde el coordinador al nodo, el nodo responde enviando una trama hacia el coordinador. El coordinador"
$line3 = "entonces compara el payload y la direccion origen para proceder a responder con otra trama."
$fullText = $line1 + [char]13 + $line2 + [char]13 + $line3 + [char]13
$shp.TextFrame.Characters().Text = $fullText
$shp.TextFrame.Characters(1, $line1.Length).Font.Bold = $true

# View changes
$ws.Activate()
$ws.Range("S18").Select()
